$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '24.632.00'
Set-TextValue $ws.Cells.Item(2, 5) '  +3.80%  '

Set-TextValue $ws.Cells.Item(3, 4) '1.695.90'
Set-TextValue $ws.Cells.Item(3, 5) '  +2.08%  '

Set-TextValue $ws.Cells.Item(4, 4) '1.002'
Set-TextValue $ws.Cells.Item(4, 5) '  -0.03%  '

Set-TextValue $ws.Cells.Item(5, 4) '318.99'
Set-TextValue $ws.Cells.Item(5, 5) '  +3.49%  '

Set-TextValue $ws.Cells.Item(6, 5) '  -0.01%  '

Set-TextValue $ws.Cells.Item(7, 4) '0.3958'
Set-TextValue $ws.Cells.Item(7, 5) '  +1.63%  '

Set-TextValue $ws.Cells.Item(8, 4) '0.4025'
Set-TextValue $ws.Cells.Item(8, 5) '  +1.77%  '

Set-TextValue $ws.Cells.Item(9, 4) '1.539'
Set-TextValue $ws.Cells.Item(9, 5) '  +8.89%  '

Set-TextValue $ws.Cells.Item(10, 2) 'BinanceUSD'
Set-TextValue $ws.Cells.Item(10, 3) 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Cells.Item(10, 4) '1.002'
Set-TextValue $ws.Cells.Item(10, 5) '  +0.02%  '

Set-TextValue $ws.Cells.Item(11, 2) 'OKB'
Set-TextValue $ws.Cells.Item(11, 3) 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Cells.Item(11, 4) '53.89'
Set-TextValue $ws.Cells.Item(11, 5) '  +8.45%  '

Set-TextValue $ws.Cells.Item(12, 4) '0.08802'
Set-TextValue $ws.Cells.Item(12, 5) '  +1.59%  '

Set-TextValue $ws.Cells.Item(13, 4) '7.263'
Set-TextValue $ws.Cells.Item(13, 5) '  +7.74%  '

Set-TextValue $ws.Cells.Item(14, 4) '23.32'
Set-TextValue $ws.Cells.Item(14, 5) '  +2.54%  '

Set-TextValue $ws.Cells.Item(15, 4) '0.00001322'
Set-TextValue $ws.Cells.Item(15, 5) '  +0.66%  '

Set-TextValue $ws.Cells.Item(16, 4) '7.630'
Set-TextValue $ws.Cells.Item(16, 5) '  +5.53%  '

Set-TextValue $ws.Cells.Item(17, 4) '1.697.17'
Set-TextValue $ws.Cells.Item(17, 5) '  +1.87%  '

Set-TextValue $ws.Cells.Item(18, 4) '101.23'
Set-TextValue $ws.Cells.Item(18, 5) '  +1.29%  '

Set-TextValue $ws.Cells.Item(19, 4) '0.07026'
Set-TextValue $ws.Cells.Item(19, 5) '  +3.69%  '

Set-TextValue $ws.Cells.Item(20, 4) '19.74'
Set-TextValue $ws.Cells.Item(20, 5) '  +3.56%  '

Set-TextValue $ws.Cells.Item(21, 4) '6.899'
Set-TextValue $ws.Cells.Item(21, 5) '  +3.59%  '

Set-TextValue $ws.Cells.Item(22, 4) '1.001'
Set-TextValue $ws.Cells.Item(22, 5) '  -0.09%  '

Set-TextValue $ws.Cells.Item(23, 4) '14.11'
Set-TextValue $ws.Cells.Item(23, 5) '  +2.19%  '

Set-TextValue $ws.Cells.Item(24, 4) '24.632.66'
Set-TextValue $ws.Cells.Item(24, 5) '  +3.85%  '

Set-TextValue $ws.Cells.Item(25, 4) '3.068'
Set-TextValue $ws.Cells.Item(25, 5) '  +9.98%  '

Set-TextValue $ws.Cells.Item(26, 4) '2.336'
Set-TextValue $ws.Cells.Item(26, 5) '  +0.92%  '

Set-TextValue $ws.Cells.Item(27, 4) '22.38'
Set-TextValue $ws.Cells.Item(27, 5) '  +2.87%  '

Set-TextValue $ws.Cells.Item(28, 4) '159.96'
Set-TextValue $ws.Cells.Item(28, 5) '  +1.60%  '

Set-TextValue $ws.Cells.Item(29, 4) '5.250'
Set-TextValue $ws.Cells.Item(29, 5) '  +1.44%  '

Set-TextValue $ws.Cells.Item(30, 4) '134.36'
Set-TextValue $ws.Cells.Item(30, 5) '  +3.70%  '

Set-TextValue $ws.Cells.Item(31, 4) '7.431'
Set-TextValue $ws.Cells.Item(31, 5) '  +15.02%  '

Set-TextValue $ws.Cells.Item(32, 4) '1.115'
Set-TextValue $ws.Cells.Item(32, 5) '  -1.45%  '

Set-TextValue $ws.Cells.Item(33, 4) '1.888.35'
Set-TextValue $ws.Cells.Item(33, 5) '  +2.12%  '

Set-TextValue $ws.Cells.Item(34, 4) '7.460'
Set-TextValue $ws.Cells.Item(34, 5) '  +15.43%  '

Set-TextValue $ws.Cells.Item(35, 4) '0.08544'
Set-TextValue $ws.Cells.Item(35, 5) '  -0.77%  '

Set-TextValue $ws.Cells.Item(36, 4) '11.47'
Set-TextValue $ws.Cells.Item(36, 5) '  +11.38%  '

Set-TextValue $ws.Cells.Item(37, 2) 'WEMIXTOKEN'
Set-TextValue $ws.Cells.Item(37, 3) 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Cells.Item(37, 4) '1.960'
Set-TextValue $ws.Cells.Item(37, 5) '  -0.94%  '

Set-TextValue $ws.Cells.Item(38, 2) 'Algorand'
Set-TextValue $ws.Cells.Item(38, 3) 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Cells.Item(38, 4) '0.2748'
Set-TextValue $ws.Cells.Item(38, 5) '  +3.68%  '

Set-TextValue $ws.Cells.Item(39, 4) '14.64'
Set-TextValue $ws.Cells.Item(39, 5) '  +1.84%  '

Set-TextValue $ws.Cells.Item(40, 4) '0.02790'
Set-TextValue $ws.Cells.Item(40, 5) '  +10.30%  '

Set-TextValue $ws.Cells.Item(41, 4) '0.09056'
Set-TextValue $ws.Cells.Item(41, 5) '  +3.24%  '

Set-TextValue $ws.Cells.Item(42, 2) 'TrustWalletToken'
Set-TextValue $ws.Cells.Item(42, 3) 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Cells.Item(42, 4) '1.468'
Set-TextValue $ws.Cells.Item(42, 5) '  +1.32%  '

Set-TextValue $ws.Cells.Item(43, 2) 'TheSandbox'
Set-TextValue $ws.Cells.Item(43, 3) 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Cells.Item(43, 4) '0.7744'
Set-TextValue $ws.Cells.Item(43, 5) '  +2.55%  '

Set-TextValue $ws.Cells.Item(44, 4) '0.7237'
Set-TextValue $ws.Cells.Item(44, 5) '  +2.97%  '

Set-TextValue $ws.Cells.Item(45, 4) '15.49'
Set-TextValue $ws.Cells.Item(45, 5) '  +4.09%  '

Set-TextValue $ws.Cells.Item(46, 4) '2.522'
Set-TextValue $ws.Cells.Item(46, 5) '  +5.94%  '

Set-TextValue $ws.Cells.Item(47, 4) '4.237'
Set-TextValue $ws.Cells.Item(47, 5) '  +3.92%  '

Set-TextValue $ws.Cells.Item(48, 4) '1.367'
Set-TextValue $ws.Cells.Item(48, 5) '  +15.94%  '

Set-TextValue $ws.Cells.Item(49, 4) '1.001'
Set-TextValue $ws.Cells.Item(49, 5) '  -0.04%  '

Set-TextValue $ws.Cells.Item(50, 4) '141.11'
Set-TextValue $ws.Cells.Item(50, 5) '  +2.15%  '

Set-TextValue $ws.Cells.Item(51, 4) '0.08048'
Set-TextValue $ws.Cells.Item(51, 5) '  +3.83%  '
